# Update the LR-pair worksheet with freshly computed TPM-based values.
# The new data set only keeps "Resolving-Mac" as the target cluster, with
# four sending clusters (ECs, FAPs, MuSCs, Resolving-Mac) instead of the
# previous 4 sending x 2 target combinations (8 data rows -> 4 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row 2..5), columns A..T
$data = @(
    @("ECs",         "C3", "Itgax", "Resolving-Mac", 3, 1, 0.555934,           1.667802,    0.005745252779589096, 0.005745252779589094, 3, 1, 44.94181766666667, 134.825453, 1, 1, 24.98468446270067, 224.862160164306,  0.005745252779589096, 0.005745252779589094),
    @("FAPs",        "C3", "Itgax", "Resolving-Mac", 3, 1, 79.68771233333334,  239.063137,  0.823525905561055,    0.823525905561055,    3, 1, 44.94181766666667, 134.825453, 1, 1, 3581.310637958452, 32231.79574162607, 0.823525905561055,    0.823525905561055),
    @("MuSCs",       "C3", "Itgax", "Resolving-Mac", 3, 1, 0.3446996666666666, 1.034099,    0.003562269474506148, 0.003562269474506148, 3, 1, 44.94181766666667, 134.825453, 1, 1, 15.49142956909411, 139.422866121847,  0.003562269474506148, 0.003562269474506148),
    @("Resolving-Mac","C3","Itgax", "Resolving-Mac", 3, 1, 16.17571666666667,  48.52715,    0.1671665721848498,   0.1671665721848498,   3, 1, 44.94181766666667, 134.825453, 1, 1, 726.9661090609945, 6542.694981548951, 0.1671665721848498,   0.1671665721848498)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

# Remove the now-unused trailing rows (previously rows 6-9).
$ws.Rows("6:9").Delete()
